$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "onerilen satis fiyati" (I) and "rakip ortalama satis fiyati" (J)
# for rows 2-4 to reflect corrected/filtered scraped product data.

$ws.Range("I2").Value = 291
$ws.Range("J2").Value = 291

$ws.Range("I3").Value = 960.8647826086956
$ws.Range("J3").Value = 960.8647826086956

$ws.Range("I4").Value = 1022.195
$ws.Range("J4").Value = 1022.195
